# Undo/Redo activity diagram: rename "address book" references to
# "financial planner" in the two shapes that mention it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "TextBox 47": "[command commits address book]" ---
$shpCommand = $s.Shapes.Item("TextBox 47")
$trCommand = $shpCommand.TextFrame.TextRange

$oldCommandRun = "command commits address book]"
$newCommandRun = "command commits financial planner]"
$commandRunStart = "[".Length + 1
$commandRun = $trCommand.Characters($commandRunStart, $oldCommandRun.Length)
$commandRun.Text = $newCommandRun

# --- Shape "Rounded Rectangle 50": "Purge redundant states..." ---
$shpPurge = $s.Shapes.Item("Rounded Rectangle 50")
$trPurge = $shpPurge.TextFrame.TextRange

$oldFirstRun = "Purge redundant states and then save address book to "
$newFirstRun = "Purge redundant states and then save financial planner to "
$firstRun = $trPurge.Characters(1, $oldFirstRun.Length)
$firstRun.Text = $newFirstRun

$oldSecondRun = "addressBookStateList"
$newSecondRun = "financialPlannerStateList"
$secondRunStart = $newFirstRun.Length + 1
$secondRun = $trPurge.Characters($secondRunStart, $oldSecondRun.Length)
$secondRun.Text = $newSecondRun
